# [Kadastro App] Yeni kayit eklendi: 2997
#
# Appends the new Erdemli kadastro record (Kayit No 2997) as row 58 to both
# the master "Kayitlar" sheet and the filtered "Erdemli" sheet, mirroring
# the existing rows (which are all plain text values, e.g. record numbers
# and dates are stored as text, not numbers/dates).

$wb = $excel.ActiveWorkbook

$newRow = @{
    A = "2997"
    B = "2025-09-11"
    C = "Erdemli"
    D = "1"
    E = "3B"
    F = "EMİNE ALANLI KIRCILI (K.Mühendisi), SEVİL SARAÇER (Tekniker)"
}

$targetSheets = @("Kayitlar", "Erdemli")

foreach ($sheetName in $targetSheets) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Force the new cells to stay text (existing rows store numbers/dates
    # as literal text, e.g. "2997" / "2025-09-11" / "1"), otherwise Excel
    # would auto-coerce these values into real numbers/dates.
    $rowRange = $ws.Range("A58:F58")
    $rowRange.NumberFormat = "@"

    $ws.Range("A58").Value = $newRow.A
    $ws.Range("B58").Value = $newRow.B
    $ws.Range("C58").Value = $newRow.C
    $ws.Range("D58").Value = $newRow.D
    $ws.Range("E58").Value = $newRow.E
    $ws.Range("F58").Value = $newRow.F
}
